# Update the math-drill worksheet to the next day's values.
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Header date
Replace-Text "2024-02-01 Thursday" "2024-02-02 Friday"

# Row 1
Replace-Text "507×9=4563" "697×2=1394"
Replace-Text "760×3=2280" "155×4=620"
# NOTE: "424×4=1696" is both an old value (cell 4) and the *new* value of
# cell 3 below, so replace the old occurrence first to avoid clobbering
# the freshly written text.
Replace-Text "424×4=1696" "213×5=1065"
Replace-Text "445×2=890" "424×4=1696"
Replace-Text "705×5=3525" "836×4=3344"

# Row 2 (problems row)
Replace-Text "103×9=927" "700×5=3500"
Replace-Text "976×4=3904" "294×5=1470"
Replace-Text "244×9=2196" "498×2=996"
Replace-Text "133×6=798" "270×7=1890"
Replace-Text "577×4=2308" "355×5=1775"

# Row 3 (problems row)
Replace-Text "459×5=2295" "731×3=2193"
Replace-Text "462×9=4158" "340×9=3060"
Replace-Text "654×6=3924" "920×6=5520"
Replace-Text "868×7=6076" "113×8=904"
Replace-Text "203×9=1827" "803×7=5621"

# Row 4 (problems row)
Replace-Text "467×4=1868" "314×7=2198"
Replace-Text "782×6=4692" "638×6=3828"
Replace-Text "446×8=3568" "778×9=7002"
Replace-Text "343×2=686" "386×4=1544"
Replace-Text "998×7=6986" "318×5=1590"

# Row 5 (problems row)
Replace-Text "466×6=2796" "726×9=6534"
Replace-Text "639×9=5751" "468×4=1872"
Replace-Text "323×4=1292" "887×6=5322"
Replace-Text "690×3=2070" "500×4=2000"
Replace-Text "814×9=7326" "481×5=2405"
